$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header of the first column from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Update the "Gen" column values (A2:A14) to the new MaxFES fractional values
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# 3. Remove the "Run 50" column (AZ) entirely - this shifts the "Mean" column
#    (previously BA) left into AZ, and shrinks the used range to A1:AZ14
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 4. Recalculate the "Mean" column (now AZ) values, now that "Run 50" data
#    has been excluded from the average
$ws.Range("AZ2").Value = 243.83729526
$ws.Range("AZ3").Value = 142.75766372
$ws.Range("AZ4").Value = 5.49966392
$ws.Range("AZ5").Value = 0.79751049
$ws.Range("AZ6").Value = 0.79751049
$ws.Range("AZ7").Value = 0.79751049
$ws.Range("AZ8").Value = 0.79751049
$ws.Range("AZ9").Value = 0.79751049
$ws.Range("AZ10").Value = 0.79751049
$ws.Range("AZ11").Value = 0.79751049
$ws.Range("AZ12").Value = 0.79751049
$ws.Range("AZ13").Value = 0.79751049
$ws.Range("AZ14").Value = 0.79751049
